# Auto-generated Excel COM-interop script
# Applies updated market-price / profit figures to the Sephirot_Profits workbook
# (mirrors a scheduled runner refreshing cached Universalis price data).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9 (Leve Item ID 5487)
$ws.Range("H9").Value = 81.333336
$ws.Range("I9").Value = 79.59999999999999
$ws.Range("K9").Value = 79.59999999999999
$ws.Range("M9").Value = 89.40000000000001

# Row 12 (Leve Item ID 5515)
$ws.Range("H12").Value = 156.33333
$ws.Range("I12").Value = 134.5
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 134.5
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = 35.5
$ws.Range("N12").Value = -540

# Row 20 (Leve Item ID 1965)
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

# Row 35 (Leve Item ID 1965)
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 5999
$ws.Range("I40").Value = 2997
$ws.Range("K40").Value = 2997
$ws.Range("M40").Value = -2822

# Row 51 (Leve Item ID 5486)
$ws.Range("H51").Value = 15000
$ws.Range("I51").Value = 15000
$ws.Range("K51").Value = 15000
$ws.Range("M51").Value = -14516

# Row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 1880
$ws.Range("I70").Value = 1880
$ws.Range("K70").Value = 5640
$ws.Range("M70").Value = -5370

# Row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 1880
$ws.Range("I73").Value = 1880
$ws.Range("K73").Value = 5640
$ws.Range("M73").Value = -4704

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 9557.6
$ws.Range("I116").Value = 9259
$ws.Range("K116").Value = 9259
$ws.Range("M116").Value = -5817

# Row 123 (Leve Item ID 34090)
$ws.Range("H123").Value = 50000
$ws.Range("I123").Value = 50000
$ws.Range("K123").Value = 50000
$ws.Range("M123").Value = -45100

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 1354.7391
$ws.Range("I132").Value = 1103.0476
$ws.Range("J132").Value = 3997.5
$ws.Range("K132").Value = 3309.142800000001
$ws.Range("L132").Value = 11992.5
$ws.Range("M132").Value = -779.1428000000005
$ws.Range("N132").Value = -17052.5

# Row 135 (Leve Item ID 44047)
$ws.Range("H135").Value = 372
$ws.Range("I135").Value = 372
$ws.Range("K135").Value = 3348
$ws.Range("M135").Value = -813

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1568.6
$ws.Range("I137").Value = 1710.75
$ws.Range("J137").Value = 1000
$ws.Range("K137").Value = 5132.25
$ws.Range("L137").Value = 3000
$ws.Range("M137").Value = -2582.25
$ws.Range("N137").Value = -8100

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 4127.409
$ws.Range("J138").Value = 4068.5789
$ws.Range("L138").Value = 12205.7367
$ws.Range("N138").Value = -22485.7367

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 4417.48
$ws.Range("I32").Value = 4518.2085
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 4518.2085
$ws.Range("L32").Value = 2000
$ws.Range("M32").Value = -4231.2085
$ws.Range("N32").Value = -2574

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 1501
$ws.Range("I74").Value = 1501
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1501
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -627
$ws.Range("N74").ClearContents()

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 1501
$ws.Range("I77").Value = 1501
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 7505
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -3137
$ws.Range("N77").ClearContents()

# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 5497.6665
$ws.Range("I97").Value = 4782.5713
$ws.Range("J97").Value = 8000.5
$ws.Range("K97").Value = 4782.5713
$ws.Range("L97").Value = 8000.5
$ws.Range("M97").Value = -4286.5713
$ws.Range("N97").Value = -8992.5

$ws = $wb.Worksheets.Item("CRP")
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 3102.5
$ws.Range("I132").Value = 1206
$ws.Range("K132").Value = 3618
$ws.Range("M132").Value = -1088

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (Leve Item ID 4650)
$ws.Range("H4").Value = 4006754.5
$ws.Range("I4").Value = 83626
$ws.Range("J4").Value = 13422263
$ws.Range("K4").Value = 250878
$ws.Range("L4").Value = 40266789
$ws.Range("M4").Value = -250766
$ws.Range("N4").Value = -40267013

# Row 114 (Leve Item ID 27865)
$ws.Range("H114").Value = 1591.5
$ws.Range("I114").Value = 1366.6666
$ws.Range("J114").Value = 1726.4
$ws.Range("K114").Value = 4099.9998
$ws.Range("L114").Value = 5179.200000000001
$ws.Range("M114").Value = -845.9997999999996
$ws.Range("N114").Value = -11687.2

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 2201.3076
$ws.Range("J131").Value = 2182.25
$ws.Range("L131").Value = 6546.75
$ws.Range("N131").Value = -16626.75

# Row 138 (Leve Item ID 44105)
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Row 107 (Leve Item ID 27802)
$ws.Range("H107").Value = 501
$ws.Range("I107").Value = 501
$ws.Range("K107").Value = 501
$ws.Range("M107").Value = 1419

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 5386.2
$ws.Range("I22").Value = 1542
$ws.Range("J22").Value = 8749.875
$ws.Range("K22").Value = 1542
$ws.Range("L22").Value = 8749.875
$ws.Range("M22").Value = -1247
$ws.Range("N22").Value = -9339.875

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 5386.2
$ws.Range("I27").Value = 1542
$ws.Range("J27").Value = 8749.875
$ws.Range("K27").Value = 1542
$ws.Range("L27").Value = 8749.875
$ws.Range("M27").Value = -1435
$ws.Range("N27").Value = -8963.875

# Row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 25540.545
$ws.Range("I82").Value = 18242.875
$ws.Range("J82").Value = 45001
$ws.Range("K82").Value = 18242.875
$ws.Range("L82").Value = 45001
$ws.Range("M82").Value = -17881.875
$ws.Range("N82").Value = -45723

# Row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 25540.545
$ws.Range("I85").Value = 18242.875
$ws.Range("J85").Value = 45001
$ws.Range("K85").Value = 18242.875
$ws.Range("L85").Value = 45001
$ws.Range("M85").Value = -16994.875
$ws.Range("N85").Value = -47497

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 3881.4443
$ws.Range("I132").Value = 3419.1428
$ws.Range("K132").Value = 10257.4284
$ws.Range("M132").Value = -7727.428400000001
